# Add a "2022" data column (T) to the mortality-rate table on the active
# sheet, mirroring the formatting already used by column S (the "2021"
# column), and move the active selection the way the author's Excel left it
# (cell U4, just to the right of the new column's header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column S's cell formatting (number format, font, borders, alignment)
# onto the new column T for the same rows (header row 4 through data row 14)
# before filling in the 2022 values.
$ws.Range("S4:S14").Copy()
$ws.Range("T4:T14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "2022" column values (one row below the existing "2021" column, S).
$ws.Range("T4").Value  = 2022   # header year
$ws.Range("T5").Value  = 3.7
$ws.Range("T6").Value  = 1.6
$ws.Range("T7").Value  = 1.7
$ws.Range("T8").Value  = 17.9
$ws.Range("T9").Value  = 7.5
$ws.Range("T10").Value = 1.1
$ws.Range("T11").Value = 4.4
$ws.Range("T12").Value = 3
$ws.Range("T13").Value = 4.1
$ws.Range("T14").Value = 0.8

# Match the workbook's saved selection/active cell (U4).
$ws.Range("U4").Select()
